$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value (2023-09-13 = 45182).
# Update every populated row (2 through 496) to the new date serial 45184
# (2023-09-15), matching the source workbook's "last checked/changed" stamp.
$newDate = [DateTime]::FromOADate(45184)

for ($row = 2; $row -le 496; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
